# Auto-generated font-fix script
# Renames font families so the 'typeface' actually matches the installed family:
#   Satoshi -> Arial
#   Ogg     -> Ogg TRIAL
$p = $ppt.ActivePresentation

$edits = @(
    @(1, 3, "Arial"),
    @(1, 4, "Ogg TRIAL"),
    @(1, 5, "Ogg TRIAL"),
    @(1, 6, "Arial"),
    @(1, 7, "Arial"),
    @(2, 4, "Ogg TRIAL"),
    @(2, 6, "Arial"),
    @(2, 7, "Arial"),
    @(3, 4, "Ogg TRIAL"),
    @(3, 8, "Arial"),
    @(4, 2, "Ogg TRIAL"),
    @(4, 5, "Arial"),
    @(4, 6, "Ogg TRIAL"),
    @(4, 9, "Arial"),
    @(4, 10, "Ogg TRIAL"),
    @(4, 11, "Ogg TRIAL"),
    @(5, 1, "Ogg TRIAL"),
    @(5, 2, "Ogg TRIAL"),
    @(5, 3, "Ogg TRIAL"),
    @(5, 4, "Arial"),
    @(6, 3, "Ogg TRIAL"),
    @(7, 1, "Ogg TRIAL"),
    @(7, 3, "Arial"),
    @(8, 3, "Ogg TRIAL"),
    @(8, 5, "Arial"),
    @(8, 7, "Arial"),
    @(9, 1, "Ogg TRIAL"),
    @(10, 1, "Ogg TRIAL"),
    @(10, 2, "Ogg TRIAL"),
    @(11, 2, "Arial"),
    @(11, 3, "Ogg TRIAL"),
    @(11, 4, "Ogg TRIAL"),
    @(12, 2, "Ogg TRIAL"),
    @(12, 4, "Arial"),
    @(12, 5, "Arial"),
    @(13, 2, "Ogg TRIAL"),
    @(13, 3, "Arial"),
    @(13, 165, "Arial"),
    @(14, 2, "Ogg TRIAL"),
    @(14, 3, "Arial"),
    @(14, 5, "Arial"),
    @(14, 9, "Arial"),
    @(14, 10, "Arial"),
    @(15, 2, "Ogg TRIAL"),
    @(15, 3, "Arial"),
    @(15, 4, "Ogg TRIAL"),
    @(15, 5, "Ogg TRIAL"),
    @(15, 6, "Ogg TRIAL"),
    @(15, 7, "Ogg TRIAL"),
    @(15, 8, "Ogg TRIAL"),
    @(15, 9, "Ogg TRIAL"),
    @(15, 10, "Ogg TRIAL"),
    @(15, 11, "Ogg TRIAL"),
    @(15, 12, "Ogg TRIAL"),
    @(15, 13, "Ogg TRIAL"),
    @(15, 14, "Ogg TRIAL"),
    @(15, 15, "Ogg TRIAL"),
    @(15, 16, "Ogg TRIAL"),
    @(15, 17, "Ogg TRIAL"),
    @(15, 18, "Ogg TRIAL"),
    @(15, 20, "Arial"),
    @(16, 2, "Ogg TRIAL"),
    @(16, 3, "Arial"),
    @(16, 4, "Arial"),
    @(17, 2, "Arial"),
    @(17, 3, "Ogg TRIAL"),
    @(17, 5, "Ogg TRIAL"),
    @(17, 6, "Ogg TRIAL"),
    @(18, 1, "Ogg TRIAL"),
    @(19, 3, "Ogg TRIAL"),
    @(19, 5, "Arial"),
    @(19, 6, "Ogg TRIAL"),
    @(20, 2, "Arial"),
    @(20, 4, "Arial"),
    @(20, 5, "Ogg TRIAL"),
    @(20, 6, "Arial"),
    @(20, 8, "Ogg TRIAL"),
    @(20, 9, "Arial"),
)

foreach ($edit in $edits) {
    $slideIndex = $edit[0]
    $shapeIndex = $edit[1]
    $fontName = $edit[2]
    $shape = $p.Slides.Item($slideIndex).Shapes.Item($shapeIndex)
    $shape.TextFrame.TextRange.Font.Name = $fontName
}
